# "Link Triggered.xlsx" — update of the "Feb" sheet's metric table
# (Links Triggered / Response / NPS / Concern Count / CC per 1000 / OSAT)
# plus the remembered cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feb")
$ws.Activate()

# ---------------------------------------------------------------------
# 1) Numeric cell updates (columns B, C, D, F) — literal values, no formulas
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 100
$ws.Range("F2").Value = 10

$ws.Range("B3").Value = 140
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 14
$ws.Range("F3").Value = 14.29

$ws.Range("B4").Value = 116
$ws.Range("F4").Value = 8.62

$ws.Range("B5").Value = 98
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 71

$ws.Range("B6").Value = 139
$ws.Range("F6").Value = 7.19

$ws.Range("B7").Value = 92

$ws.Range("B8").Value = 120

$ws.Range("B9").Value = 112
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 50
$ws.Range("F9").Value = 8.93

$ws.Range("B12").Value = 109

$ws.Range("B14").Value = 67
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 0
$ws.Range("F14").Value = 14.93

$ws.Range("B15").Value = 112
$ws.Range("F15").Value = 8.93

$ws.Range("B17").Value = 98
$ws.Range("F17").Value = 10.2

$ws.Range("B18").Value = 111
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 57

$ws.Range("B19").Value = 34
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 57

$ws.Range("B20").Value = 94
$ws.Range("F20").Value = 10.64

$ws.Range("B22").Value = 65
$ws.Range("F22").Value = 15.38

$ws.Range("B23").Value = 85
$ws.Range("F23").Value = 11.76

$ws.Range("B24").Value = 123

$ws.Range("B25").Value = 123

$ws.Range("B26").Value = 60
$ws.Range("F26").Value = 16.67

$ws.Range("B27").Value = 138
$ws.Range("F27").Value = 28.99

$ws.Range("B28").Value = 94
$ws.Range("F28").Value = 21.28

$ws.Range("B29").Value = 72

$ws.Range("B30").Value = 3

$ws.Range("B31").Value = 97

$ws.Range("B32").Value = 82
$ws.Range("F32").Value = 24.39

$ws.Range("B33").Value = 50
$ws.Range("F33").Value = 20

$ws.Range("B34").Value = 48
$ws.Range("F34").Value = 20.83

$ws.Range("B35").Value = 56
$ws.Range("F35").Value = 17.86

# ---------------------------------------------------------------------
# 2) OSAT (column G) updates — these are literal percentage TEXT labels
#    (same shared-string style already used throughout the workbook), not
#    numbers. Typing "43%" directly would make Excel coerce it into a
#    percent-formatted number, so instead copy the already-existing text
#    label from elsewhere in the workbook (values-only paste keeps the
#    destination cell's own formatting untouched).
# ---------------------------------------------------------------------
$wbSheets = $wb.Worksheets

$wbSheets.Item("Apr").Range("G29").Copy()
$ws.Range("G3").PasteSpecial(-4163)   # "43%"

$wbSheets.Item("Apr").Range("G8").Copy()
$ws.Range("G5").PasteSpecial(-4163)   # "71%"

$ws.Range("G22").Copy()
$ws.Range("G9").PasteSpecial(-4163)   # "60%"

$wbSheets.Item("Nov").Range("G18").Copy()
$ws.Range("G18").PasteSpecial(-4163)  # "29%"

$wbSheets.Item("May").Range("G26").Copy()
$ws.Range("G19").PasteSpecial(-4163)  # "43%"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Column widths (B:G) — the sheet previously only carried explicit
#    widths for columns A, H and I; bring B:G in line with the rest of
#    the workbook's sheets (best-effort — this runtime's column-width
#    engine rounds to whole display pixels, so exact fractional widths
#    from a real-Excel "best fit" aren't reproducible bit-for-bit).
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 12.4094
$ws.Columns.Item(3).ColumnWidth = 8.0773
$ws.Columns.Item(4).ColumnWidth = 3.9041
$ws.Columns.Item(5).ColumnWidth = 12.3987
$ws.Columns.Item(6).ColumnWidth = 7.24
$ws.Columns.Item(7).ColumnWidth = 4.7414

# ---------------------------------------------------------------------
# 4) Selection moves from H12 to H4
# ---------------------------------------------------------------------
$ws.Range("H4").Select()
